# "Solucionada parte fechas de autor"
# - AUTOR table (columns K/L on "Modelo físico"): the two DATE columns
#   (fecha_nacimiento / fecha_fallecimiento) become VARCHAR(4) (year only).
# - LIBRO table (columns H/I on "Modelo físico"): the rows were reordered so
#   that the FK columns (autor [FK] / coautor [FK]) come right after the PK,
#   followed by titulo / paginas.
# - Selection/scroll position of that sheet is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "Modelo físico"

# --- AUTOR table: DATE -> VARCHAR(4) -----------------------------------
$ws.Range("L9").Value  = "VARCHAR(4)"
$ws.Range("L10").Value = "VARCHAR(4)"

# --- LIBRO table: reorder H7:I10 ----------------------------------------
# Before:                       After:
#  H7  titulo      / VARCHAR NOT NULL       H7  autor [FK]   / INT REFERENCES autor
#  H8  paginas     / SMALLINT NOT NULL      H8  coautor [FK] / INT REFERENCES autor
#  H9  autor [FK]  / INT REFERENCES autor   H9  titulo       / VARCHAR NOT NULL
#  H10 coautor [FK]/ INT REFERENCES autor   H10 paginas      / SMALLINT NOT NULL
$ws.Range("H7").Value  = "autor [FK]"
$ws.Range("I7").Value  = "INT REFERENCES autor"
$ws.Range("H8").Value  = "coautor [FK]"
$ws.Range("I8").Value  = "INT REFERENCES autor"
$ws.Range("H9").Value  = "titulo"
$ws.Range("I9").Value  = "VARCHAR NOT NULL"
$ws.Range("H10").Value = "paginas"
$ws.Range("I10").Value = "SMALLINT NOT NULL"

# --- Update the view: scrolled to column C, selection on L11 -----------
$ws.Activate()
$ws.Range("L11").Select()
$excel.ActiveWindow.ScrollColumn = 3
